$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update total_customers (C8) and returning_customers (E8) for year 2025
$ws.Range("C8").Value = 853
$ws.Range("E8").Value = 708

# Recompute new_rate (G8) and returning_rate (H8) based on updated totals
$ws.Range("G8").Value = 83.00117233294256
$ws.Range("H8").Value = 16.99882766705744
